# Update metadata workbook with 2p parameters for preprocessing.
# End state: three sheets - animals, stacks, two_photon_settings.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# --- Sheet 1: rename existing "stacks" sheet to "animals" and rebuild its header row ---
$wsAnimals = $wb.Worksheets.Item(1)
$wsAnimals.Name = "animals"

$animalsHeaders = @("animal_id", "genotype", "owner", "tank", "date_of_birth", "root_dir")
for ($i = 0; $i -lt $animalsHeaders.Count; $i++) {
    $wsAnimals.Cells.Item(1, $i + 1).Value = $animalsHeaders[$i]
}
# Drop the old trailing columns (original sheet went out to X1) so the used range shrinks back to F1.
$wsAnimals.Range("G1:X1").Clear()

# Keep a reference to a still-styled header cell (style index "1": bold/bordered/centered) to
# stamp onto the brand-new sheets below via PasteSpecial(Formats), since freshly-added sheets
# start with no cell styling at all.
$styleSource = $wsAnimals.Range("A1")

# --- Sheet 2: new "stacks" sheet (rebuilt with microscope_settings_path / num_planes removed) ---
$wsStacks = $wb.Worksheets.Add($null, $wsAnimals)
$wsStacks.Name = "stacks"

$stacksHeaders = @(
    "animal_id",
    "stack_id",
    "stack_type",
    "date",
    "condition",
    "experimenter",
    "include_in_analysis",
    "image_quality",
    "notes",
    "raw_path",
    "stimulus_name",
    "stimulus_metadata_path",
    "zoom_factor",
    "round_id",
    "plane_spacing",
    "reference_channel_index",
    "channel1_name",
    "channel1_wavelength_nm",
    "channel2_name",
    "channel2_wavelength_nm",
    "channel3_name",
    "channel3_wavelength_nm"
)
for ($i = 0; $i -lt $stacksHeaders.Count; $i++) {
    $wsStacks.Cells.Item(1, $i + 1).Value = $stacksHeaders[$i]
}
$wsStacks.PageSetup.LeftMargin = 54
$wsStacks.PageSetup.RightMargin = 54
$wsStacks.PageSetup.TopMargin = 72
$wsStacks.PageSetup.BottomMargin = 72
$wsStacks.PageSetup.HeaderMargin = 36
$wsStacks.PageSetup.FooterMargin = 36

# --- Sheet 3: new "two_photon_settings" sheet ---
$wsTwoPhoton = $wb.Worksheets.Add($null, $wsStacks)
$wsTwoPhoton.Name = "two_photon_settings"

$twoPhotonHeaders = @(
    "session_id",
    "mode",
    "n_planes",
    "frames_per_plane",
    "flyback_frames",
    "remove_first_frame",
    "blocks"
)
for ($i = 0; $i -lt $twoPhotonHeaders.Count; $i++) {
    $wsTwoPhoton.Cells.Item(1, $i + 1).Value = $twoPhotonHeaders[$i]
}
$wsTwoPhoton.PageSetup.LeftMargin = 54
$wsTwoPhoton.PageSetup.RightMargin = 54
$wsTwoPhoton.PageSetup.TopMargin = 72
$wsTwoPhoton.PageSetup.BottomMargin = 72
$wsTwoPhoton.PageSetup.HeaderMargin = 36
$wsTwoPhoton.PageSetup.FooterMargin = 36

# Stamp the bold/bordered/centered header style onto the new sheets' header rows
# (values already populated above, so a formats-only paste just restyles them).
$styleSource.Copy()
$wsStacks.Range("A1:V1").PasteSpecial($xlPasteFormats)
$wsTwoPhoton.Range("A1:G1").PasteSpecial($xlPasteFormats)

# Make sure the first sheet ("animals") is the active one, matching activeTab="0".
$wsAnimals.Activate() | Out-Null
$wsAnimals.Range("A1").Select() | Out-Null
